$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the scores for row 23 (student "Сидаков Амир")
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = 5

# Update the active selection to F23, matching the latest view state
$ws.Range("F23").Select()
